$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "VHF transmitter not used meaning there is a hole in the bottom of the drop-off plate."
$ws.Range("B7").Value = "NA"
$ws.Range("C7").Value = "Fill the hole with silicon sealant and allow it to cure completely. This will technically be a 'weak' point in the drop-off plate. However, the location of this point, i.e., against the belting, means it will face significant exposure to the environment once the SensorDrop system detaches. "
